$d = $word.ActiveDocument

# Merge the runs of the "Poisson Probability Distribution..." TOC paragraph
# (paragraph 21) into a single run while preserving its original formatting.
# Trick: inserting a character right at the boundary of the first run causes
# Word to reflow/merge adjacent runs that share identical formatting; we then
# delete the inserted placeholder character, leaving one merged run.
$pPoisson = $d.Paragraphs.Item(21)
$rPoisson = $pPoisson.Range
$rPoisson.End = $rPoisson.End - 1
$poissonRun1Start = $rPoisson.Start
$poissonRun1End = $poissonRun1Start + 7
$poissonRun1 = $d.Range($poissonRun1Start, $poissonRun1End)
$poissonRun1.InsertAfter("X")
$poissonPlaceholder = $d.Range($poissonRun1End, $poissonRun1End + 1)
$poissonPlaceholder.Delete()

# Merge the runs of the "Chebyshev's Theorem..." TOC paragraph (paragraph 22)
# into a single run the same way.
$pCheb = $d.Paragraphs.Item(22)
$rCheb = $pCheb.Range
$rCheb.End = $rCheb.End - 1
$chebRun1Start = $rCheb.Start
$chebRun1End = $chebRun1Start + 19
$chebRun1 = $d.Range($chebRun1Start, $chebRun1End)
$chebRun1.InsertAfter("X")
$chebPlaceholder = $d.Range($chebRun1End, $chebRun1End + 1)
$chebPlaceholder.Delete()

# Right-align the ten Table of Contents entry paragraphs.
for ($i = 20; $i -le 29; $i++) {
    $entry = $d.Paragraphs.Item($i)
    $entry.Alignment = 2
}
